# Add a new "Deepseek_error" column (L) to the LLM-eval results sheet:
#   L1 -> header "Deepseek_error"
#   L2 -> the (wrapped) SPARQL query text for that run
# Then restore the selection/view state to match where the author ended up.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column header -----------------------------------------------------
$ws.Range("L1").Value = "Deepseek_error"

# --- New column value (row 2), wrapped like the rest of the data column ---
$deepseekError = "PREFIX bfo: <http://purl.obolibrary.org/obo/bfo.owl/>`nPREFIX cdio: <https://w3id.org/CDIO/>`nPREFIX dc: <http://purl.org/dc/elements/1.1/>`nPREFIX ns1: <http://purl.obolibrary.org/obo/bfo.owl#>`nPREFIX obi: <http://purl.obolibrary.org/obo/obi.owl/>`nPREFIX xsd: <http://www.w3.org/2001/XMLSchema#>`nSELECT DISTINCT ?study`nWHERE {`n  GRAPH <https://w3id.org/CDIO/graph/studies> {`n    ?study a obi:StudyDesignExecution .`n    ?study bfo:concretizes ?study_design .`n    ?study_design bfo:has_part ?protocol .`n    # Look for diabetes in study primary purpose`n    {`n      ?protocol bfo:has_part ?purpose_spec .`n      ?purpose_spec ?has_value ?purpose_value .`n      FILTER (CONTAINS(LCASE(STR(?purpose_value)), `"diabetes`"))`n    }`n    UNION`n    # Optionally look for diabetes in inclusion criteria`n    {`n      ?protocol bfo:has_part ?inclusion_spec .`n      ?inclusion_spec ?has_value ?inclusion_value .`n      FILTER (CONTAINS(LCASE(STR(?inclusion_value)), `"diabetes`"))`n    }`n  }`n}"

$ws.Range("L2").Value = $deepseekError
$ws.Range("L2").WrapText = $true

# --- Column width for the new column (matches column C's width) -----------
$ws.Columns.Item(12).ColumnWidth = 16.5

# --- Selection / view state -------------------------------------------------
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 2
[void]$ws.Range("O2").Select()
